$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '67.456.81'
$ws.Range('E2').Value = '  -0.21%  '

# Row 3
$ws.Range('D3').Value = '2.521.32'
$ws.Range('E3').Value = '  +0.20%  '

# Row 4
$ws.Range('E4').Value = '  +0.04%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '588.92'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.45%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '170.07'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.60%  '

# Row 7
$ws.Range('E7').Value = '  +0.07%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.521'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.89%  '

# Row 9
$ws.Range('D9').Value = '2.520.95'
$ws.Range('E9').Value = '  +0.22%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.136'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.32%  '

# Row 11
$ws.Range('E11').Value = '  +1.06%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.341'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.92%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.97'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -3.69%  '

# Row 14
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '2.988.98'
$ws.Range('E14').Value = '  +0.41%  '

# Row 15
$ws.Range('B15').Value = 'Avalanche'
$ws.Range('C15').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '26.04'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.71%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000174'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.42%  '

# Row 17
$ws.Range('D17').Value = '67.578.56'
$ws.Range('E17').Value = '  +0.21%  '

# Row 18
$ws.Range('D18').Value = '2.540.70'
$ws.Range('E18').Value = '  +0.93%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '8.03'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.15%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.66'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.28%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '363.48'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.52%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.14'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.77%  '

# Row 23
$ws.Range('B23').Value = 'NEARProtocol'
$ws.Range('C23').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.49'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.17%  '

# Row 24
$ws.Range('B24').Value = 'Litecoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '71.76'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.05%  '

# Row 25
$ws.Range('B25').Value = 'Dai'
$ws.Range('C25').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.00'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.01%  '

# Row 26
$ws.Range('B26').Value = 'SuiNetwork'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.84'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -6.13%  '

# Row 27
$ws.Range('B27').Value = 'Aptos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.72'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -4.92%  '

# Row 28
$ws.Range('B28').Value = 'WrappedeETH'
$ws.Range('C28').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D28').Value = '2.662.35'
$ws.Range('E28').Value = '  +0.71%  '

# Row 29
$ws.Range('B29').Value = 'PEPE'
$ws.Range('C29').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D29').Value = '0.0₃0940'
$ws.Range('E29').Value = '  -4.60%  '

# Row 30
$ws.Range('B30').Value = 'Bittensor'
$ws.Range('C30').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '528.94'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.62%  '

# Row 31
$ws.Range('B31').Value = 'InternetComputer(DFINITY)'
$ws.Range('C31').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.20'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.30%  '

# Row 32
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.85'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.68%  '

# Row 33
$ws.Range('B33').Value = 'Fetch.AI'
$ws.Range('C33').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.28'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -4.73%  '

# Row 34
$ws.Range('B34').Value = 'FirstDigitalUSD'
$ws.Range('C34').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.00'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.00%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.127'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -2.22%  '

# Row 36
$ws.Range('B36').Value = 'Monero'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '157.66'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.33%  '

# Row 37
$ws.Range('B37').Value = 'EthereumClassic'
$ws.Range('C37').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '19.33'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.90%  '

# Row 38
$ws.Range('B38').Value = 'ImmutableX'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.42'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.62%  '

# Row 39
$ws.Range('B39').Value = 'WhiteBITCoin'
$ws.Range('C39').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.62'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.11%  '

# Row 40
$ws.Range('B40').Value = 'Stacks'
$ws.Range('C40').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.76'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.53%  '

# Row 41
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '5.02'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -2.96%  '

# Row 42
$ws.Range('B42').Value = 'PolygonEcosystemToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.341'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.14%  '

# Row 43
$ws.Range('B43').Value = 'USDe'
$ws.Range('C43').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.00'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.18%  '

# Row 44
$ws.Range('B44').Value = 'dogwifhat'
$ws.Range('C44').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.42'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -3.73%  '

# Row 45
$ws.Range('B45').Value = 'OKB'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '39.36'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.37%  '

# Row 46
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '147.26'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.73%  '

# Row 47
$ws.Range('B47').Value = 'Filecoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.68'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.81%  '

# Row 48
$ws.Range('B48').Value = 'ARBITRUM'
$ws.Range('C48').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.545'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.64%  '

# Row 49
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').Value = '0.0…0271'
$ws.Range('E49').Value = '  -2.43%  '

# Row 50
$ws.Range('B50').Value = 'Optimism'
$ws.Range('C50').Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.69'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.58%  '

# Row 51
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.595'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.37%  '
